$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (col G) values for rows 2-68, computed from the
# updated Strike/s_vals calculation.
$kVals = @(
    0, 2, 0, 1, 0, 1, 0, 1, 1, 2, 0, 0, 2, 1, 2, 2, 3, 1, 1, 2, 1, 2, 0, 1, 4, 0, 1, 0, 0, 1, 0, 0, 1, 0, 1, 0, 0, 1, 2, 1, 1, 2, 2, 2, 4, 0, 1, 1, 2, 1, 2, 2, 0, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 2, 2, 1, 0
)

for ($i = 0; $i -lt $kVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kVals[$i]
}

